# Natmi following Dr Hou advice
# Recomputes the LR-pairs table (Sema4a -> Plxnb1) for all 4x4 sending/target
# cluster combinations (ECs, FAPs, M2, sCs), expanding the sheet from 8 data
# rows (A1:T9) to 16 data rows (A1:T17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema4a"
$ws.Range("C2").Value = "Plxnb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 22.614608
$ws.Range("H2").Value = 67.843824
$ws.Range("I2").Value = 0.3650188533124966
$ws.Range("J2").Value = 0.3650188533124966
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.7376096666666667
$ws.Range("N2").Value = 2.212829
$ws.Range("O2").Value = 0.1014284037978316
$ws.Range("P2").Value = 0.1014284037978316
$ws.Range("Q2").Value = 16.68075346867733
$ws.Range("R2").Value = 150.126781218096
$ws.Range("S2").Value = 0.03702327964760136
$ws.Range("T2").Value = 0.03702327964760135

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema4a"
$ws.Range("C3").Value = "Plxnb1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 22.614608
$ws.Range("H3").Value = 67.843824
$ws.Range("I3").Value = 0.3650188533124966
$ws.Range("J3").Value = 0.3650188533124966
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.6291593333333333
$ws.Range("N3").Value = 1.887478
$ws.Range("O3").Value = 0.08651544278546762
$ws.Range("P3").Value = 0.08651544278546762
$ws.Range("Q3").Value = 14.22819169287467
$ws.Range("R3").Value = 128.053725235872
$ws.Range("S3").Value = 0.0315797677193743
$ws.Range("T3").Value = 0.03157976771937429

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema4a"
$ws.Range("C4").Value = "Plxnb1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 22.614608
$ws.Range("H4").Value = 67.843824
$ws.Range("I4").Value = 0.3650188533124966
$ws.Range("J4").Value = 0.3650188533124966
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.07810833333333334
$ws.Range("N4").Value = 0.234325
$ws.Range("O4").Value = 0.01074064499332162
$ws.Range("P4").Value = 0.01074064499332162
$ws.Range("Q4").Value = 1.766389339866667
$ws.Range("R4").Value = 15.8975040588
$ws.Range("S4").Value = 0.003920537919298865
$ws.Range("T4").Value = 0.003920537919298864

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Sema4a"
$ws.Range("C5").Value = "Plxnb1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 22.614608
$ws.Range("H5").Value = 67.843824
$ws.Range("I5").Value = 0.3650188533124966
$ws.Range("J5").Value = 0.3650188533124966
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.827342666666667
$ws.Range("N5").Value = 17.482028
$ws.Range("O5").Value = 0.8013155084233792
$ws.Range("P5").Value = 0.8013155084233792
$ws.Range("Q5").Value = 131.7830700883413
$ws.Range("R5").Value = 1186.047630795072
$ws.Range("S5").Value = 0.2924952680262221
$ws.Range("T5").Value = 0.2924952680262221

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sema4a"
$ws.Range("C6").Value = "Plxnb1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 12.38193366666667
$ws.Range("H6").Value = 37.145801
$ws.Range("I6").Value = 0.1998548561530699
$ws.Range("J6").Value = 0.1998548561530699
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.7376096666666667
$ws.Range("N6").Value = 2.212829
$ws.Range("O6").Value = 0.1014284037978316
$ws.Range("P6").Value = 0.1014284037978316
$ws.Range("Q6").Value = 9.133033964558779
$ws.Range("R6").Value = 82.19730568102901
$ws.Range("S6").Value = 0.02027095905085112
$ws.Range("T6").Value = 0.02027095905085112

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema4a"
$ws.Range("C7").Value = "Plxnb1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 12.38193366666667
$ws.Range("H7").Value = 37.145801
$ws.Range("I7").Value = 0.1998548561530699
$ws.Range("J7").Value = 0.1998548561530699
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.6291593333333333
$ws.Range("N7").Value = 1.887478
$ws.Range("O7").Value = 0.08651544278546762
$ws.Range("P7").Value = 0.08651544278546762
$ws.Range("Q7").Value = 7.790209131097555
$ws.Range("R7").Value = 70.111882179878
$ws.Range("S7").Value = 0.01729053137290878
$ws.Range("T7").Value = 0.01729053137290878

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Sema4a"
$ws.Range("C8").Value = "Plxnb1"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 12.38193366666667
$ws.Range("H8").Value = 37.145801
$ws.Range("I8").Value = 0.1998548561530699
$ws.Range("J8").Value = 0.1998548561530699
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.07810833333333334
$ws.Range("N8").Value = 0.234325
$ws.Range("O8").Value = 0.01074064499332162
$ws.Range("P8").Value = 0.01074064499332162
$ws.Range("Q8").Value = 0.9671322021472223
$ws.Range("R8").Value = 8.704189819325
$ws.Range("S8").Value = 0.002146570060131482
$ws.Range("T8").Value = 0.002146570060131482

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Sema4a"
$ws.Range("C9").Value = "Plxnb1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 12.38193366666667
$ws.Range("H9").Value = 37.145801
$ws.Range("I9").Value = 0.1998548561530699
$ws.Range("J9").Value = 0.1998548561530699
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 5.827342666666667
$ws.Range("N9").Value = 17.482028
$ws.Range("O9").Value = 0.8013155084233792
$ws.Range("P9").Value = 0.8013155084233792
$ws.Range("Q9").Value = 72.15377035160311
$ws.Range("R9").Value = 649.383933164428
$ws.Range("S9").Value = 0.1601467956691785
$ws.Range("T9").Value = 0.1601467956691785

$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Sema4a"
$ws.Range("C10").Value = "Plxnb1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 19.82277733333333
$ws.Range("H10").Value = 59.468332
$ws.Range("I10").Value = 0.3199563508543806
$ws.Range("J10").Value = 0.3199563508543806
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.7376096666666667
$ws.Range("N10").Value = 2.212829
$ws.Range("O10").Value = 0.1014284037978316
$ws.Range("P10").Value = 0.1014284037978316
$ws.Range("Q10").Value = 14.62147218124755
$ws.Range("R10").Value = 131.593249631228
$ws.Range("S10").Value = 0.03245266195213879
$ws.Range("T10").Value = 0.03245266195213879

$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Sema4a"
$ws.Range("C11").Value = "Plxnb1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 19.82277733333333
$ws.Range("H11").Value = 59.468332
$ws.Range("I11").Value = 0.3199563508543806
$ws.Range("J11").Value = 0.3199563508543806
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.6291593333333333
$ws.Range("N11").Value = 1.887478
$ws.Range("O11").Value = 0.08651544278546762
$ws.Range("P11").Value = 0.08651544278546762
$ws.Range("Q11").Value = 12.47168537185511
$ws.Range("R11").Value = 112.245168346696
$ws.Range("S11").Value = 0.02768116536618917
$ws.Range("T11").Value = 0.02768116536618917

$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Sema4a"
$ws.Range("C12").Value = "Plxnb1"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 19.82277733333333
$ws.Range("H12").Value = 59.468332
$ws.Range("I12").Value = 0.3199563508543806
$ws.Range("J12").Value = 0.3199563508543806
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.07810833333333334
$ws.Range("N12").Value = 0.234325
$ws.Range("O12").Value = 0.01074064499332162
$ws.Range("P12").Value = 0.01074064499332162
$ws.Range("Q12").Value = 1.548324099544444
$ws.Range("R12").Value = 13.9349168959
$ws.Range("S12").Value = 0.003436537577885558
$ws.Range("T12").Value = 0.003436537577885558

$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Sema4a"
$ws.Range("C13").Value = "Plxnb1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 19.82277733333333
$ws.Range("H13").Value = 59.468332
$ws.Range("I13").Value = 0.3199563508543806
$ws.Range("J13").Value = 0.3199563508543806
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 5.827342666666667
$ws.Range("N13").Value = 17.482028
$ws.Range("O13").Value = 0.8013155084233792
$ws.Range("P13").Value = 0.8013155084233792
$ws.Range("Q13").Value = 115.5141161263662
$ws.Range("R13").Value = 1039.627045137296
$ws.Range("S13").Value = 0.256385985958167
$ws.Range("T13").Value = 0.256385985958167

$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Sema4a"
$ws.Range("C14").Value = "Plxnb1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 7.135311000000001
$ws.Range("H14").Value = 21.405933
$ws.Range("I14").Value = 0.115169939680053
$ws.Range("J14").Value = 0.115169939680053
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.7376096666666667
$ws.Range("N14").Value = 2.212829
$ws.Range("O14").Value = 0.1014284037978316
$ws.Range("P14").Value = 0.1014284037978316
$ws.Range("Q14").Value = 5.263074368273001
$ws.Range("R14").Value = 47.367669314457
$ws.Range("S14").Value = 0.01168150314724032
$ws.Range("T14").Value = 0.01168150314724032

$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Sema4a"
$ws.Range("C15").Value = "Plxnb1"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 7.135311000000001
$ws.Range("H15").Value = 21.405933
$ws.Range("I15").Value = 0.115169939680053
$ws.Range("J15").Value = 0.115169939680053
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.6291593333333333
$ws.Range("N15").Value = 1.887478
$ws.Range("O15").Value = 0.08651544278546762
$ws.Range("P15").Value = 0.08651544278546762
$ws.Range("Q15").Value = 4.489247511886
$ws.Range("R15").Value = 40.403227606974
$ws.Range("S15").Value = 0.00996397832699538
$ws.Range("T15").Value = 0.00996397832699538

$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Sema4a"
$ws.Range("C16").Value = "Plxnb1"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 7.135311000000001
$ws.Range("H16").Value = 21.405933
$ws.Range("I16").Value = 0.115169939680053
$ws.Range("J16").Value = 0.115169939680053
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.07810833333333334
$ws.Range("N16").Value = 0.234325
$ws.Range("O16").Value = 0.01074064499332162
$ws.Range("P16").Value = 0.01074064499332162
$ws.Range("Q16").Value = 0.557327250025
$ws.Range("R16").Value = 5.015945250225
$ws.Range("S16").Value = 0.001236999436005714
$ws.Range("T16").Value = 0.001236999436005714

$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Sema4a"
$ws.Range("C17").Value = "Plxnb1"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 7.135311000000001
$ws.Range("H17").Value = 21.405933
$ws.Range("I17").Value = 0.115169939680053
$ws.Range("J17").Value = 0.115169939680053
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 5.827342666666667
$ws.Range("N17").Value = 17.482028
$ws.Range("O17").Value = 0.8013155084233792
$ws.Range("P17").Value = 0.8013155084233792
$ws.Range("Q17").Value = 41.579902230236
$ws.Range("R17").Value = 374.219120072124
$ws.Range("S17").Value = 0.09228745876981156
$ws.Range("T17").Value = 0.09228745876981156

